# "Made changes to Benefits in Policy Modal"
# Add two new header cells (Z1/AA1) with a light-blue fill, used as a
# "Help_text" / "help_desc" column pair appended to the existing param grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("Z1").Value = "Help_text"
$ws.Range("AA1").Value = "help_desc"

# Highlight the new header cells with a solid light-blue fill (RGB 00B0F0)
$ws.Range("Z1:AA1").Interior.Color = 15773696

# Match the author's final selection/view state
$ws.Range("Z1:AA1").Select() | Out-Null
